$d = $word.ActiveDocument

function Insert-ItalicParagraphAfter($searchText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $endPos = $rng.End
    $rng.InsertParagraphAfter()
    $newRng = $d.Range($endPos + 1, $endPos + 1)
    $newRng.InsertAfter($newText)
    $newRng.Italic = 1
}

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# 2. Insert English translation paragraph after "Objetivos" (Portuguese) text
Insert-ItalicParagraphAfter "Apresentar uma visão geral da química dos elementos e de seus compostos enfatizando as correlações entre as propriedades físicas e químicas com os aspectos estruturais e de ligação, os métodos de obtenção em laboratório e indústria, além das principais propriedades e aplicações." "To present an overview of the chemistry of elements and their compounds, emphasizing the correlations between physical and chemical properties with structural and binding aspects, methods of obtaining them in laboratory and industry, in addition to the main properties and applications."

# 3. Insert English translation paragraph after "Programa resumido" (Portuguese) text
Insert-ItalicParagraphAfter "Ocorrência, obtenção, estrutura, propriedades e aplicações de elementos metálicos e não-metálicos; moléculas poliatômicas; compostos halogenados e das famílias do oxigênio, nitrogênio, carbono e boro; compostos oxigenados. Processos industriais de fabricação." "Occurrence, obtaining, structure, properties and applications of metallic and non-metallic elements; polyatomic molecules; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds. Industrial manufacturing processes."

# 4. Insert English translation paragraph after "Programa" (Portuguese) text
Insert-ItalicParagraphAfter "Ocorrência, obtenção, estrutura, propriedades e aplicações de elementos não-metálicos: gases nobres, hidrogênio molecular, halogênios, oxigênio molecular, ozônio e nitrogênio molecular; semimetais; metais alcalinos, alcalinos-terrosos e de transição; moléculas poliatômicas e espécies catenadas de: enxofre, fósforo e carbono; compostos halogenados e das famílias do oxigênio, nitrogênio, carbono e boro; compostos oxigenados: óxidos e oxicompostos. Processos industriais de fabricação dos principais insumos químicos e materiais." "Occurrence, obtaining, structure, properties and applications of non-metallic elements: noble gases, molecular hydrogen, halogens, molecular oxygen, ozone and molecular nitrogen; semimetals; alkali, alkaline earth and transition metals; polyatomic molecules and catenated species of: sulfur, phosphorus and carbon; halogenated compounds and the oxygen, nitrogen, carbon and boron families; oxygenated compounds: oxides and oxycompounds. Industrial manufacturing processes of the main chemical inputs and materials."

# 5. Update requisite course text
$d.Content.Find.Execute("LOQ4031 -  Química Geral I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)", 2)
